$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'knee pads snowboard'
$ws.Range("A2").Value = 'youth under armour basketball tights'
$ws.Range("A3").Value = 'black basketball rim'
$ws.Range("A4").Value = 'mcdavid youth knee pads basketball'
$ws.Range("A5").Value = 'nike compression pants men'
$ws.Range("A6").Value = 'asics knee pad'
$ws.Range("A7").Value = 'nike leggings basketball'
$ws.Range("A8").Value = 'compression tights with knee pads'
$ws.Range("A9").Value = 'combat knee pads'
$ws.Range("A10").Value = 'uflex knee compression'
$ws.Range("A11").Value = 'damascus knee pads'
$ws.Range("A12").Value = 'mcdavid youth basketball knee pads'
$ws.Range("A13").Value = 'mens pants black'
$ws.Range("A14").Value = 'basketball pants girls'
$ws.Range("A15").Value = 'compression capri leggings'
$ws.Range("A16").Value = 'padded tights'
$ws.Range("A17").Value = 'padded tights men'
$ws.Range("A18").Value = 'crye knee pads'
$ws.Range("A19").Value = 'dancing knee pads'
$ws.Range("A20").Value = 'dodgeball knee pads'
$ws.Range("A21").Value = 'uflex athletics knee compression sleeve'
$ws.Range("A22").Value = 'gform knee pads'
$ws.Range("A23").Value = 'nike tights for men'
$ws.Range("A24").Value = 'insertable knee pads'
$ws.Range("A25").Value = 'knee pads for scootering'
$ws.Range("A26").Value = 'youth compression pants with pads'
$ws.Range("A27").Value = 'red knee pads'
$ws.Range("A28").Value = 'black volleyball knee pads'
$ws.Range("A29").Value = 'men basketball pants'
$ws.Range("A30").Value = 'compression tights with pads'
$ws.Range("A31").Value = 'compression leggings with knee pads'
$ws.Range("A32").Value = 'ski knee pads'
$ws.Range("A33").Value = 'dead on knee pads'
$ws.Range("A34").Value = 'black pads'
$ws.Range("A35").Value = 'padded basketball compression pants'
$ws.Range("A36").Value = 'mens padded leggings'
$ws.Range("A37").Value = 'knee pads for skating'
$ws.Range("A38").Value = 'padded compression tights'
$ws.Range("A39").Value = 'basketball knee pads pants'
$ws.Range("A40").Value = 'compression pants womens'
$ws.Range("A41").Value = 'military knee pad'
$ws.Range("A42").Value = 'drskin mens compression pants'
$ws.Range("A43").Value = 'padded compression pants men basketball'
$ws.Range("A44").Value = 'snowboard knee pad'
$ws.Range("A45").Value = 'double knee pads'
$ws.Range("A46").Value = 'knee pads light'
$ws.Range("A47").Value = 'knee pad protection'
$ws.Range("A48").Value = 'leggings with knee pads women'
$ws.Range("A49").Value = 'padded tights men basketball'
$ws.Range("A50").Value = 'mens leggings with knee pads'
$ws.Range("A51").Value = 'knee pad under pants'
$ws.Range("A52").Value = 'football pants adult xxl'
$ws.Range("A53").Value = 'the best leggings'
$ws.Range("A54").Value = 'sliding shorts mens'
$ws.Range("A55").Value = 'snowboarding padded shorts men'
$ws.Range("A56").Value = 'flexible knee pads'
$ws.Range("A57").Value = 'basketball shorts for men pack'
$ws.Range("A58").Value = 'leggings capri'
$ws.Range("A59").Value = 'leggings spandex'
$ws.Range("A60").Value = 'floor knee pads'
$ws.Range("A61").Value = 'gym pads'
$ws.Range("A62").Value = 'pants with knee pads kids'
$ws.Range("A63").Value = 'basketball pants men'
$ws.Range("A64").Value = 'volleyball knee pads nike black'
$ws.Range("A65").Value = 'tights men'
$ws.Range("A66").Value = 'hayabusa compression pants'
$ws.Range("A67").Value = 'tough knee pads'
$ws.Range("A68").Value = 'knee pad leggings'
$ws.Range("A69").Value = 'goalkeeper pant'
$ws.Range("A70").Value = 'knee compression men'
$ws.Range("A71").Value = 'wrestling knee sleeves'
$ws.Range("A72").Value = 'padded knee sleeves for basketball'
$ws.Range("A73").Value = 'knee compression sleeve for squats'
$ws.Range("A74").Value = 'youth leg sleeves for basketball'
$ws.Range("A75").Value = 'tactical knee pads'
$ws.Range("A76").Value = 'knee pads paintball'
$ws.Range("A77").Value = 'bike knee pads'
$ws.Range("A78").Value = 'knee pads biking'
$ws.Range("A79").Value = 'goalie knee pads'
$ws.Range("A80").Value = 'knee pads compression sleeve'
$ws.Range("A81").Value = 'pant with knee pads'
$ws.Range("A82").Value = 'long knee pads'
$ws.Range("A83").Value = 'knee pads for girls'
$ws.Range("A84").Value = 'mens pants with knee pads'
$ws.Range("A85").Value = 'maroon knee pads'
$ws.Range("A86").Value = 'pantalones con rodilleras'
$ws.Range("A87").Value = 'pantalon con rodilleras'
$ws.Range("A88").Value = 'knee padded pants men'
$ws.Range("A89").Value = 'knee pad pants men'
$ws.Range("A90").Value = 'baseball sliding pants mens'
$ws.Range("A91").Value = 'men basketball knee pads'
$ws.Range("A92").Value = 'baseball pants mens knee'
$ws.Range("A93").Value = 'compression leggings basketball'
$ws.Range("A94").Value = 'black mens baseball pants'
$ws.Range("A95").Value = 'basketball pants for men'
$ws.Range("A96").Value = 'work pants with knee pads'
$ws.Range("A97").Value = 'knee pad for basketball youth'
$ws.Range("A98").Value = 'basketball tights for men'
$ws.Range("A99").Value = 'black athletic pants men'
$ws.Range("A100").Value = 'youth basketball clothes'
